$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the note text in K2
$ws.Range("K2").Value = "*Note Type => General / Birthday"

# Bold, blue (Accent1 theme color) font for the note
$ws.Range("K2").Font.Bold = $true
$ws.Range("K2").Font.ThemeColor = 5   # msoThemeColorAccent1

# Column K width to match the diff (~29 characters)
$ws.Columns.Item(11).ColumnWidth = 29

# Update selection to match the diff (K10 becomes the active cell)
$ws.Range("K10").Select()
